$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''62.265.24'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.38%  '
$ws.Range("D3").Value = '''2.422.86'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.78%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '''564.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.18%  '
$ws.Range("D6").Value = '''144.72'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.42%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +1.84%  '
$ws.Range("D9").Value = '''2.421.06'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.62%  '
$ws.Range("E10").Value = '  +1.80%  '
$ws.Range("D12").Value = '''5.40'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.02%  '
$ws.Range("E13").Value = '  +0.86%  '
$ws.Range("D14").Value = '''26.07'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.15%  '
$ws.Range("D15").Value = '''0.0000177'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.27%  '
$ws.Range("D16").Value = '''2.861.23'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.81%  '
$ws.Range("D17").Value = '''62.112.35'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.28%  '
$ws.Range("D18").Value = '''2.419.57'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.73%  '
$ws.Range("D19").Value = '''11.31'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.61%  '
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '''4.20'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.06%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = '''325.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.24%  '
$ws.Range("E22").Value = '  +0.40%  '
$ws.Range("D24").Value = '''65.62'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.92%  '
$ws.Range("E25").Value = '  -2.95%  '
$ws.Range("D26").Value = '''9.03'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.26%  '
$ws.Range("D27").Value = '''590.64'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +13.94%  '
$ws.Range("D28").Value = '''2.545.21'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.91%  '
$ws.Range("D29").Value = '''0.0₃0948'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.00%  '
$ws.Range("E30").Value = '  -0.23%  '
$ws.Range("D31").Value = '''1.47'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.35%  '
$ws.Range("D32").Value = '''8.28'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.79%  '
$ws.Range("D33").Value = '''0.150'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.24%  '
$ws.Range("D34").Value = '''1.88'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.01%  '
$ws.Range("E35").Value = '  +0.86%  '
$ws.Range("E36").Value = '  +3.85%  '
$ws.Range("D37").Value = '''1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D38").Value = '''4.80'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.16%  '
$ws.Range("D39").Value = '''153.69'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.19%  '
$ws.Range("E40").Value = '  +1.12%  '
$ws.Range("E41").Value = '  +1.01%  '
$ws.Range("D42").Value = '''1.83'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.64%  '
$ws.Range("D43").Value = '''0.999'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("D44").Value = '''2.35'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.61%  '
$ws.Range("D45").Value = '''150.21'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.17%  '
$ws.Range("E46").Value = '  +1.46%  '
$ws.Range("E47").Value = '  +2.33%  '
$ws.Range("D48").Value = '''20.43'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.49%  '
$ws.Range("E49").Value = '  +1.92%  '
$ws.Range("E50").Value = '  +2.20%  '
$ws.Range("E51").Value = '  +1.86%  '
